$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.6893
$ws.Range("E2").Value = 16.32140000000001
$ws.Range("B3").Value = 5.863
$ws.Range("D3").Value = -7.216099999999996
$ws.Range("E6").Value = 16.45430000000001
$ws.Range("D12").Value = -7.397099999999999
$ws.Range("B14").Value = 5.768100000000003
$ws.Range("B16").Value = 6.225400000000004
$ws.Range("C18").Value = -12.1274
$ws.Range("E19").Value = 16.1923
$ws.Range("B21").Value = 8.849700000000006
$ws.Range("B23").Value = 8.611700000000006
$ws.Range("C24").Value = -12.6383
$ws.Range("D24").Value = -8.099899999999996
$ws.Range("E24").Value = 16.52990000000001
$ws.Range("B25").Value = 5.389700000000004
$ws.Range("C25").Value = -13.1181
$ws.Range("D25").Value = -9.056899999999992
$ws.Range("B26").Value = 6.626399999999999
$ws.Range("C27").Value = -11.48789999999999
$ws.Range("E27").Value = 16.95049999999999
$ws.Range("B29").Value = 5.131700000000004
$ws.Range("C30").Value = -13.2856
$ws.Range("E30").Value = 16.041
$ws.Range("C31").Value = -13.63499999999999
$ws.Range("E31").Value = 16.13890000000001
$ws.Range("E33").Value = 17.05790000000001
$ws.Range("C39").Value = -12.73730000000001
$ws.Range("B40").Value = 9.186299999999999
$ws.Range("D41").Value = -7.852399999999996
$ws.Range("C42").Value = -12.4347
$ws.Range("E42").Value = 16.24750000000001
$ws.Range("C48").Value = -11.59339999999999
$ws.Range("D50").Value = -8.214700000000001
$ws.Range("C51").Value = -11.49899999999999
$ws.Range("C52").Value = -11.4996
$ws.Range("B53").Value = 5.179100000000001
$ws.Range("D53").Value = -5.688500000000001
$ws.Range("C55").Value = -13.3831
$ws.Range("E55").Value = 16.51620000000001
$ws.Range("C56").Value = -11.2408
$ws.Range("D56").Value = -8.153700000000001
$ws.Range("B57").Value = 5.173599999999999
$ws.Range("C57").Value = -13.7434
$ws.Range("D57").Value = -8.939399999999997
$ws.Range("D58").Value = -8.519800000000007
$ws.Range("E58").Value = 16.22200000000002
$ws.Range("B59").Value = 4.409499999999998
$ws.Range("C60").Value = -13.57279999999999
$ws.Range("D61").Value = -8.136299999999999
$ws.Range("D63").Value = -7.893400000000002
$ws.Range("D64").Value = -7.916099999999997
$ws.Range("B65").Value = 5.721900000000005
$ws.Range("E65").Value = 17.16310000000001
$ws.Range("B69").Value = 5.514499999999994
$ws.Range("D70").Value = -8.1867
$ws.Range("E70").Value = 16.88949999999999
$ws.Range("D72").Value = -7.336399999999994
$ws.Range("C73").Value = -12.4115
$ws.Range("C74").Value = -12.92040000000001
$ws.Range("E74").Value = 17.12689999999998
$ws.Range("E75").Value = 16.2726
$ws.Range("B79").Value = 9.063400000000003
$ws.Range("B83").Value = 5.826499999999997
$ws.Range("E83").Value = 16.68659999999999
$ws.Range("E84").Value = 16.8172
$ws.Range("D86").Value = -8.987099999999998
$ws.Range("E86").Value = 16.05890000000001
$ws.Range("C89").Value = -10.72840000000001
$ws.Range("D89").Value = -5.633200000000004
$ws.Range("C90").Value = -12.4146
$ws.Range("B91").Value = 4.814299999999998
$ws.Range("C92").Value = -11.7086
$ws.Range("B93").Value = 5.679399999999998
$ws.Range("E96").Value = 16.7396
$ws.Range("E97").Value = 16.72670000000001
$ws.Range("D98").Value = -8.983699999999997
$ws.Range("B100").Value = 4.481299999999999
$ws.Range("D100").Value = -8.8371
$ws.Range("D102").Value = -7.678899999999999
